$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 68000
$ws.Range("J109").Value = 68000
$ws.Range("L109").Value = 68000
$ws.Range("N109").Value = -70774
$ws.Range("H112").Value = 1616.5135
$ws.Range("J112").Value = 1654.6285
$ws.Range("L112").Value = 4963.8855
$ws.Range("N112").Value = -7179.8855
$ws.Range("H138").Value = 2096.1978
$ws.Range("I138").Value = 2324.8286
$ws.Range("J138").Value = 1953.3036
$ws.Range("K138").Value = 6974.485799999999
$ws.Range("L138").Value = 5859.9108
$ws.Range("M138").Value = -1834.485799999999
$ws.Range("N138").Value = -16139.9108
$ws.Range("H139").Value = 52949.9
$ws.Range("J139").Value = 52949.9
$ws.Range("L139").Value = 52949.9
$ws.Range("N139").Value = -63229.9
$ws.Range("H140").Value = 84266.664
$ws.Range("J140").Value = 84266.664
$ws.Range("L140").Value = 84266.664
$ws.Range("N140").Value = -94626.664

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 618025.25
$ws.Range("I2").Value = 618025.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 618025.25
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -617912.25
$ws.Range("H32").Value = 6753.7334
$ws.Range("I32").Value = 4794.477
$ws.Range("K32").Value = 4794.477
$ws.Range("M32").Value = -4507.477
$ws.Range("H45").Value = 1135
$ws.Range("I45").Value = 889.5714
$ws.Range("K45").Value = 889.5714
$ws.Range("M45").Value = -512.5714
$ws.Range("H97").Value = 1637.2
$ws.Range("I97").Value = 1641.8667
$ws.Range("J97").Value = 1623.2
$ws.Range("K97").Value = 1641.8667
$ws.Range("L97").Value = 1623.2
$ws.Range("M97").Value = -1145.8667
$ws.Range("N97").Value = -2615.2
$ws.Range("H116").Value = 618025.25
$ws.Range("I116").Value = 618025.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 618025.25
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -615731.25
$ws.Range("H122").Value = 1279.9375
$ws.Range("I122").Value = 1069.9286
$ws.Range("K122").Value = 3209.7858
$ws.Range("M122").Value = -759.7857999999997
$ws.Range("H130").Value = 55820
$ws.Range("J130").Value = 55820
$ws.Range("L130").Value = 55820
$ws.Range("N130").Value = -65860
$ws.Range("H139").Value = 44000
$ws.Range("J139").Value = 44000
$ws.Range("L139").Value = 44000
$ws.Range("N139").Value = -54280

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 618025.25
$ws.Range("I3").Value = 618025.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 618025.25
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -617911.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1593.15
$ws.Range("I31").Value = 1093.1818
$ws.Range("J31").Value = 2204.2222
$ws.Range("K31").Value = 1093.1818
$ws.Range("L31").Value = 2204.2222
$ws.Range("M31").Value = -798.1818000000001
$ws.Range("N31").Value = -2794.2222
$ws.Range("H34").Value = 1593.15
$ws.Range("I34").Value = 1093.1818
$ws.Range("J34").Value = 2204.2222
$ws.Range("K34").Value = 1093.1818
$ws.Range("L34").Value = 2204.2222
$ws.Range("M34").Value = -891.1818000000001
$ws.Range("N34").Value = -2608.2222
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
$ws.Range("H135").Value = 44709
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 590.5417
$ws.Range("I5").Value = 499.8889
$ws.Range("K5").Value = 1499.6667
$ws.Range("M5").Value = -1387.6667
$ws.Range("H131").Value = 18231.525
$ws.Range("J131").Value = 18672.59
$ws.Range("L131").Value = 56017.77
$ws.Range("N131").Value = -66097.77
$ws.Range("H135").Value = 590.5417
$ws.Range("I135").Value = 499.8889
$ws.Range("K135").Value = 4499.0001
$ws.Range("M135").Value = -1964.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 67033
$ws.Range("J52").Value = 67033
$ws.Range("L52").Value = 67033
$ws.Range("N52").Value = -67551
$ws.Range("H101").Value = 3000
$ws.Range("J101").Value = 3000
$ws.Range("L101").Value = 3000
$ws.Range("N101").Value = -9490
$ws.Range("H102").Value = 2076.5715
$ws.Range("I102").Value = 2129.8823
$ws.Range("K102").Value = 2129.8823
$ws.Range("M102").Value = -507.8823000000002
$ws.Range("H107").Value = 350
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 150
$ws.Range("K107").Value = 550
$ws.Range("L107").Value = 150
$ws.Range("M107").Value = 1370
$ws.Range("N107").Value = -3990
$ws.Range("H113").Value = 1161.1428
$ws.Range("I113").Value = 970.2
$ws.Range("J113").Value = 1267.2222
$ws.Range("K113").Value = 970.2
$ws.Range("L113").Value = 1267.2222
$ws.Range("M113").Value = 1199.8
$ws.Range("N113").Value = -5607.2222
$ws.Range("H122").Value = 1639.4
$ws.Range("I122").Value = 1132.3334
$ws.Range("K122").Value = 3397.0002
$ws.Range("M122").Value = -947.0002
$ws.Range("H132").Value = 1926000.4
$ws.Range("I132").Value = 2406076.5
$ws.Range("K132").Value = 7218229.5
$ws.Range("M132").Value = -7215699.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4372.8125
$ws.Range("I16").Value = 6174.636
$ws.Range("J16").Value = 408.8
$ws.Range("K16").Value = 6174.636
$ws.Range("L16").Value = 408.8
$ws.Range("M16").Value = -6004.636
$ws.Range("N16").Value = -748.8
$ws.Range("H61").Value = 3916.25
$ws.Range("J61").Value = 4900
$ws.Range("L61").Value = 4900
$ws.Range("N61").Value = -5304
$ws.Range("H113").Value = 3916.25
$ws.Range("J113").Value = 4900
$ws.Range("L113").Value = 4900
$ws.Range("N113").Value = -9240

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1133.909
$ws.Range("I100").Value = 950.5
$ws.Range("K100").Value = 1901
$ws.Range("M100").Value = -1360
$ws.Range("H113").Value = 634.5
$ws.Range("I113").Value = 426.57144
$ws.Range("K113").Value = 1279.71432
$ws.Range("M113").Value = 890.28568
